# Applies the "Optuna Attempt (go back with original)" edit to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 1
$ws1.Range("H2").Value = 0.77
$ws1.Range("J2").Value = "Urgent"
$ws1.Range("L2").Value = 0.97

$ws1.Range("D3").Value = 1
$ws1.Range("H3").Value = 0
$ws1.Range("L3").Value = 1.02

$ws1.Range("D4").Value = 1
$ws1.Range("L4").Value = 0.89

$ws1.Range("D5").Value = 1
$ws1.Range("L5").Value = 0.88

$ws1.Range("L6").Value = 0.99
$ws1.Range("L7").Value = 1.02
$ws1.Range("L8").Value = 0.88
$ws1.Range("L9").Value = 0.96
$ws1.Range("L10").Value = 0.87
$ws1.Range("L11").Value = 1.09
$ws1.Range("L12").Value = 0.92

$ws1.Range("D13").Value = 1
$ws1.Range("L13").Value = 0.93

$ws1.Range("L14").Value = 0.95
$ws1.Range("L15").Value = 1.18
$ws1.Range("L16").Value = 0.92
$ws1.Range("L17").Value = 1.03

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "14"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "8"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "5"
